# Adds a calculated "insert" column to the Tabela1 table and appends four
# new rows of NLQ/SQL sample data to the pfp worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# --- 1. Append four new data rows (19-22) to the table -------------------
$newRowsData = @(
    @(
        "Which field produces the most oil per month?",
        "Which field produces the oil month?",
        "SELECT field, oil_production, year, month FROM ANP ",
        "SELECT year, month, field, max(oil_production) as max_oil_production FROM NLIDB_RESULT_SET GROUP BY year, month, field ORDER BY year, month, field"
    ),
    @(
        "Which basin has the highest yearly oil production?",
        "Which basin has the year oil production?",
        "SELECT basin, year, oil_production FROM ANP ",
        "SELECT basin, max(oil_production) as max_oil_production FROM (SELECT basin, SUM(oil_production) as oil_production FROM NLIDB_RESULT_SET GROUP BY basin, year) GROUP BY basin ORDER BY basin"
    ),
    @(
        "Which federated state has the lowest gas production?",
        "Which federated state has the gas production?",
        "SELECT state, gas_production FROM ANP ",
        "SELECT state, min(gas_production) as min_gas_production FROM NLIDB_RESULT_SET GROUP BY state ORDER BY state"
    ),
    @(
        "Which state of the federation has the lowest gas production?",
        "Which state of the federation has the gas production?",
        "SELECT state, gas_production FROM ANP ",
        "SELECT state, min(gas_production) as min_gas_production FROM NLIDB_RESULT_SET GROUP BY state ORDER BY state"
    )
)

foreach ($rowData in $newRowsData) {
    $newRow = $lo.ListRows.Add()
    $r = $newRow.Range
    $r.Item(1).Value = $rowData[0]
    $r.Item(2).Value = $rowData[1]
    $r.Item(3).Value = $rowData[2]
    $r.Item(4).Value = $rowData[3]
}

# --- 2. Add the calculated "insert" column --------------------------------
$newCol = $lo.ListColumns.Add()
$newCol.Range.Item(1).Value = "insert"

$formula = "=""INSERT INTO NLIDB_SQL_FROM_NLQ  VALUES('""&Tabela1[[#This Row],[NLQ preproccessed by GLAMORISE]]&""', '""&SUBSTITUTE(Tabela1[[#This Row],[NLIDB SQL]],""'"",""''"")&""');"""
$newCol.DataBodyRange.Formula = $formula

# --- 3. Cosmetics: autofit columns and update the view ---------------------
$ws.Columns.Item(1).EntireColumn.AutoFit()
$ws.Columns.Item(2).EntireColumn.AutoFit()
$ws.Columns.Item(3).EntireColumn.AutoFit()
$ws.Columns.Item(4).EntireColumn.AutoFit()
$ws.Columns.Item(5).EntireColumn.AutoFit()

$ws.Range("E19:E22").Select()
